$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date"

$ws.Range("A2").Value = 43704
$ws.Range("B2").Value = "Tues"

$ws.Range("A3").Value = 43706
$ws.Range("B3").Value = "Thurs"
$ws.Range("D3").Value = "- ``Command Line Basics <command_line_part1.ipynb>``_`n- Read and sign syllabus`n- ``Register with IPUMS <https://uma.pop.umn.edu/usa/user/new>``_`n- Register for DataCamp`n"

$ws.Range("A4").Value = 43711
$ws.Range("B4").Value = "Tues"

$ws.Range("A5").Value = 43713
$ws.Range("B5").Value = "Thurs"

$ws.Range("A6").Value = 43718
$ws.Range("B6").Value = "Tues"

$ws.Range("A7").Value = 43720
$ws.Range("B7").Value = "Thurs"

$ws.Range("A8").Value = 43725
$ws.Range("B8").Value = "Tues"

$ws.Range("A9").Value = 43727
$ws.Range("B9").Value = "Thurs"

$ws.Range("A10").Value = 43732
$ws.Range("B10").Value = "Tues"

$ws.Range("A11").Value = 43734
$ws.Range("B11").Value = "Thurs"

$ws.Range("A12").Value = 43739
$ws.Range("B12").Value = "Tues"

$ws.Range("A13").Value = 43741
$ws.Range("B13").Value = "Thurs"

$ws.Range("A14").Value = 43746
$ws.Range("B14").Value = "Tues"

$ws.Range("A15").Value = 43748
$ws.Range("B15").Value = "Thurs"

$ws.Range("A16").Value = 43753
$ws.Range("B16").Value = "Tues"

$ws.Range("A17").Value = 43755
$ws.Range("B17").Value = "Thurs"

$ws.Range("A18").Value = 43760
$ws.Range("B18").Value = "Tues"

$ws.Range("A19").Value = 43762
$ws.Range("B19").Value = "Thurs"

$ws.Range("A20").Value = 43767
$ws.Range("B20").Value = "Tues"

$ws.Range("A21").Value = 43769
$ws.Range("B21").Value = "Thurs"

$ws.Range("A22").Value = 43774
$ws.Range("B22").Value = "Tues"

$ws.Range("A23").Value = 43776
$ws.Range("B23").Value = "Thurs"

$ws.Range("A24").Value = 43781
$ws.Range("B24").Value = "Tues"

$ws.Range("A25").Value = 43783
$ws.Range("B25").Value = "Thurs"

$ws.Range("A26").Value = 43788
$ws.Range("B26").Value = "Tues"

$ws.Range("A27").Value = 43790
$ws.Range("B27").Value = "Thurs"

$ws.Range("A28").Value = 43795
$ws.Range("B28").Value = "Tues"

$ws.Range("A29").Value = 43797
$ws.Range("B29").Value = "Thurs"

$ws.Range("A30").Value = 43802
$ws.Range("B30").Value = "Tues"

$ws.Range("A31").Value = 43804
$ws.Range("B31").Value = "Thurs"

$ws.Range("A32").Value = 43809
$ws.Range("B32").Value = "Tues"

$ws.Range("A33").Value = 43811

$ws.Range("D14").Select() | Out-Null